$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# The table "Tableau1" (B4:L41) had an empty trailing row 40 (only the
# calculated "Durée" formula was present, with a blank result). The author
# filled in that row with a new journal entry, and added a new row 41 below
# it that keeps the same calculated-column formula but stays otherwise
# blank (already there before, untouched).

# Date (column B) - stored as serial 44266 => 2021-03-11. Copy the date
# style down from the row above (instead of assigning a DateTime / a fresh
# NumberFormat string) so the existing date-format cell style gets reused
# rather than a brand new one being created.
$ws.Range("B39").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B40").Value = 44266

# Heure début (C) / Heure fin (D) - time-only values (fractions of a day)
$ws.Range("C40").Value = 0.47916666666666669
$ws.Range("D40").Value = 0.5

# Module / Type / Tâche / Lieu
$ws.Range("F40").Value = "I-431"
$ws.Range("G40").Value = "Code"
$ws.Range("H40").Value = "Jeu"
$ws.Range("I40").Value = "CPNV"

# Descriptif (new entry text from the commit "Fait qu'on ne peut pas retoucher une case")
$ws.Range("J40").Value = "J'ai fait en sorte qu'on puisse pas toucher 2 fois une casse"

# Terminer
$ws.Range("K40").Value = "Oui"

# Row 40 grew taller to fit the now-wrapped Descriptif text, matching the
# other multi-line rows in the journal (28.8 = 2 x default row height).
$ws.Rows.Item(40).RowHeight = 28.8

# Update the selection so it matches where the author ended up after editing.
$ws.Range("H37").Select()
